# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Only column G ("K") values change; every other column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$gValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 3
    15 = 2
    16 = 4
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 3
    27 = 3
    28 = 1
    29 = 1
    31 = 2
    32 = 0
    33 = 3
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 2
    40 = 0
    41 = 1
    42 = 1
    43 = 3
    44 = 4
    45 = 3
    46 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
